$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": right-answer marks changed from 5 to 4, wrong-answer penalty from -1 to -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": total marks obtained changed from 60 to 48, and the "X / Y" summary updated accordingly
$ws.Range("B12").Value = 48
$ws.Range("E12").Value = "48 / 112"
